$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 815.8182
$ws.Range("I2").Value = 475.83334
$ws.Range("K2").Value = 475.83334
$ws.Range("M2").Value = -362.83334
$ws.Range("H38").Value = 362.7143
$ws.Range("I38").Value = 362.7143
$ws.Range("K38").Value = 1088.1429
$ws.Range("M38").Value = -716.1428999999998
$ws.Range("H40").Value = 1630.3334
$ws.Range("I40").Value = 1630.3334
$ws.Range("K40").Value = 1630.3334
$ws.Range("M40").Value = -1455.3334
$ws.Range("H70").Value = 100102200
$ws.Range("I70").Value = 3499.2
$ws.Range("J70").Value = 200200900
$ws.Range("K70").Value = 10497.6
$ws.Range("L70").Value = 600602700
$ws.Range("M70").Value = -10227.6
$ws.Range("N70").Value = -600603240
$ws.Range("H73").Value = 100102200
$ws.Range("I73").Value = 3499.2
$ws.Range("J73").Value = 200200900
$ws.Range("K73").Value = 10497.6
$ws.Range("L73").Value = 600602700
$ws.Range("M73").Value = -9561.599999999999
$ws.Range("N73").Value = -600604572
$ws.Range("H96").Value = 1111
$ws.Range("I96").Value = 249
$ws.Range("J96").Value = 1326.5
$ws.Range("K96").Value = 747
$ws.Range("L96").Value = 3979.5
$ws.Range("M96").Value = 626
$ws.Range("N96").Value = -6725.5
$ws.Range("H129").Value = 189462.47
$ws.Range("I129").Value = 1257.6
$ws.Range("J129").Value = 1600999
$ws.Range("K129").Value = 3772.8
$ws.Range("L129").Value = 4802997
$ws.Range("M129").Value = 1227.2
$ws.Range("N129").Value = -4812997
$ws.Range("H138").Value = 2892.2334
$ws.Range("I138").Value = 1341.8125
$ws.Range("J138").Value = 3456.0227
$ws.Range("K138").Value = 4025.4375
$ws.Range("L138").Value = 10368.0681
$ws.Range("M138").Value = 1114.5625
$ws.Range("N138").Value = -20648.0681

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1593.3478
$ws.Range("I2").Value = 1268.9524
$ws.Range("K2").Value = 1268.9524
$ws.Range("M2").Value = -1155.9524
$ws.Range("H61").Value = 6290.4287
$ws.Range("I61").Value = 4888.4116
$ws.Range("J61").Value = 12249
$ws.Range("K61").Value = 4888.4116
$ws.Range("L61").Value = 12249
$ws.Range("M61").Value = -4676.4116
$ws.Range("N61").Value = -12673
$ws.Range("H102").Value = 35421.6
$ws.Range("I102").Value = 31464.5
$ws.Range("J102").Value = 51250
$ws.Range("K102").Value = 31464.5
$ws.Range("L102").Value = 51250
$ws.Range("M102").Value = -29842.5
$ws.Range("N102").Value = -54494
$ws.Range("H110").Value = 2623.182
$ws.Range("H116").Value = 1593.3478
$ws.Range("I116").Value = 1268.9524
$ws.Range("K116").Value = 1268.9524
$ws.Range("M116").Value = 1025.0476
$ws.Range("H122").Value = 4822.936
$ws.Range("I122").Value = 4357.15
$ws.Range("K122").Value = 13071.45
$ws.Range("M122").Value = -10621.45
$ws.Range("H136").Value = 6290.4287
$ws.Range("I136").Value = 4888.4116
$ws.Range("J136").Value = 12249
$ws.Range("K136").Value = 14665.2348
$ws.Range("L136").Value = 36747
$ws.Range("M136").Value = -12115.2348
$ws.Range("N136").Value = -41847

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1593.3478
$ws.Range("I3").Value = 1268.9524
$ws.Range("K3").Value = 1268.9524
$ws.Range("M3").Value = -1154.9524
$ws.Range("H107").Value = 3107.1482
$ws.Range("I107").Value = 3405.353
$ws.Range("J107").Value = 2600.2
$ws.Range("K107").Value = 3405.353
$ws.Range("L107").Value = 2600.2
$ws.Range("M107").Value = -1485.353
$ws.Range("N107").Value = -6440.2
$ws.Range("H134").Value = 3381.8076
$ws.Range("I134").Value = 3356.1904
$ws.Range("J134").Value = 3489.4
$ws.Range("K134").Value = 10068.5712
$ws.Range("L134").Value = 10468.2
$ws.Range("M134").Value = -7533.5712
$ws.Range("N134").Value = -15538.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 12318
$ws.Range("I5").Value = 8446.666999999999
$ws.Range("K5").Value = 8446.666999999999
$ws.Range("M5").Value = -8334.666999999999
$ws.Range("H45").Value = 23067
$ws.Range("I45").Value = 23067
$ws.Range("K45").Value = 23067
$ws.Range("M45").Value = -22474
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H122").Value = 103626.53
$ws.Range("I122").Value = 140118
$ws.Range("K122").Value = 420354
$ws.Range("M122").Value = -417904
$ws.Range("H131").Value = 89298.5
$ws.Range("J131").Value = 88665
$ws.Range("L131").Value = 88665
$ws.Range("N131").Value = -98745
$ws.Range("H140").Value = 61893
$ws.Range("J140").Value = 61893
$ws.Range("L140").Value = 61893
$ws.Range("N140").Value = -72253

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 249
$ws.Range("I23").Value = 55.333332
$ws.Range("J23").Value = 301.81818
$ws.Range("K23").Value = 165.999996
$ws.Range("L23").Value = 905.45454
$ws.Range("M23").Value = 69.00000399999999
$ws.Range("N23").Value = -1375.45454
$ws.Range("H121").Value = 1000265.44
$ws.Range("I121").Value = 286
$ws.Range("K121").Value = 858
$ws.Range("M121").Value = 452
$ws.Range("H131").Value = 2811.3333
$ws.Range("I131").Value = 2159.2
$ws.Range("K131").Value = 6477.599999999999
$ws.Range("M131").Value = -1437.599999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 91611.75
$ws.Range("I12").Value = 103429.43
$ws.Range("K12").Value = 103429.43
$ws.Range("M12").Value = -103289.43
$ws.Range("H97").Value = 593.4231
$ws.Range("I97").Value = 582.17645
$ws.Range("J97").Value = 614.6667
$ws.Range("K97").Value = 582.17645
$ws.Range("L97").Value = 614.6667
$ws.Range("M97").Value = -86.17645000000005
$ws.Range("N97").Value = -1606.6667
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H102").Value = 2111.077
$ws.Range("I102").Value = 1613.762
$ws.Range("K102").Value = 1613.762
$ws.Range("M102").Value = 8.238000000000056
$ws.Range("H122").Value = 7693.4287
$ws.Range("I122").Value = 7212
$ws.Range("J122").Value = 8335.333000000001
$ws.Range("K122").Value = 21636
$ws.Range("L122").Value = 25005.999
$ws.Range("M122").Value = -19186
$ws.Range("N122").Value = -29905.999
$ws.Range("H134").Value = 8830.666999999999
$ws.Range("J134").Value = 8830.666999999999
$ws.Range("L134").Value = 26492.001
$ws.Range("N134").Value = -31562.001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3872.4546
$ws.Range("I7").Value = 3619.6
$ws.Range("J7").Value = 4083.1667
$ws.Range("K7").Value = 3619.6
$ws.Range("L7").Value = 4083.1667
$ws.Range("M7").Value = -3507.6
$ws.Range("N7").Value = -4307.1667
$ws.Range("H16").Value = 818.1739
$ws.Range("I16").Value = 373
$ws.Range("K16").Value = 373
$ws.Range("M16").Value = -203
$ws.Range("H32").Value = 16333
$ws.Range("I32").Value = 16333
$ws.Range("K32").Value = 16333
$ws.Range("M32").Value = -16016
$ws.Range("H61").Value = 2419.2307
$ws.Range("I61").Value = 1871.7778
$ws.Range("K61").Value = 1871.7778
$ws.Range("M61").Value = -1669.7778
$ws.Range("H113").Value = 2419.2307
$ws.Range("I113").Value = 1871.7778
$ws.Range("K113").Value = 1871.7778
$ws.Range("M113").Value = 298.2221999999999
$ws.Range("H118").Value = 39500
$ws.Range("J118").Value = 39000
$ws.Range("L118").Value = 39000
$ws.Range("N118").Value = -42314
$ws.Range("H122").Value = 3750
$ws.Range("J122").Value = 4375
$ws.Range("L122").Value = 13125
$ws.Range("N122").Value = -18025
$ws.Range("H126").Value = 3872.4546
$ws.Range("I126").Value = 3619.6
$ws.Range("J126").Value = 4083.1667
$ws.Range("K126").Value = 10858.8
$ws.Range("L126").Value = 12249.5001
$ws.Range("M126").Value = -8388.799999999999
$ws.Range("N126").Value = -17189.5001
$ws.Range("H136").Value = 4214.1113
$ws.Range("I136").Value = 4030.2666
$ws.Range("K136").Value = 12090.7998
$ws.Range("M136").Value = -9540.799800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 136.15
$ws.Range("J107").Value = 70.2
$ws.Range("L107").Value = 210.6
$ws.Range("N107").Value = -4050.6
$ws.Range("H122").Value = 3779.7368
$ws.Range("I122").Value = 3782.1333
$ws.Range("K122").Value = 11346.3999
$ws.Range("M122").Value = -8896.3999
$ws.Range("H136").Value = 6875.85
$ws.Range("I136").Value = 7438.857
$ws.Range("J136").Value = 5562.1665
$ws.Range("K136").Value = 22316.571
$ws.Range("L136").Value = 16686.4995
$ws.Range("N136").Value = -21786.4995
